$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 11 (nba_prediction); rows below shift up by one.
$ws.Rows.Item(11).Delete()

# Clear contents+format for columns that should end up fully blank (no style)
$ws.Range("A10").Clear()
$ws.Range("C10").Clear()
$ws.Range("D10").Clear()
$ws.Range("F10").Clear()
$ws.Range("G10").Clear()
$ws.Range("H10").Clear()
$ws.Range("I10").Clear()

# Clear contents only for B10 and E10, preserving their number format / hyperlink style
$ws.Range("B10").ClearContents()
$ws.Range("E10").ClearContents()

$ws.Range("A10:XFD10").Select()
